$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header D1 typo: "Type.1" -> "Type."
$ws.Range("D1").Value = "Type."

# Correct C12 variety name
$ws.Range("C12").Value = "Yellow Butterwax"

# Correct E21 quantity
$ws.Range("E21").Value = 20

# Add new rows 25-27
$ws.Range("A25").Value = "Veg"
$ws.Range("B25").Value = "Bean"
$ws.Range("C25").Value = "Orca"
$ws.Range("E25").Value = 24
$ws.Range("F25").Value = 3.25

$ws.Range("A26").Value = "Veg"
$ws.Range("B26").Value = "Watermelon"
$ws.Range("C26").Value = "Sugar Baby"
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 8.5

$ws.Range("A27").Value = "Veg"
$ws.Range("B27").Value = "Peas"
$ws.Range("C27").Value = "Sweet Magnolia"
$ws.Range("E27").Value = 25
$ws.Range("F27").Value = 2
